# Report params config
# - Update the Saiku REST URL stored in _settings!B1 to the new
#   rill-analysis-web REST endpoint.
# - Widen column B on _settings so the longer URL text still fits.
# - Make _settings the active/selected sheet instead of _input.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("_settings")

# New REST endpoint for the hyperlink / text shown in B1.
$settings.Range("B1").Value = "http://10.81.21.140:8280/rill-analysis-web/rest/"

# Column B needs to be widened to fit the longer URL text.
$settings.Columns.Item(2).ColumnWidth = 53.142857142857

# Switch the selected/active sheet from _input to _settings.
$settings.Activate()
